$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.417.27'
$ws.Range("E2").Value = '  +3.41%  '
$ws.Range("D3").Value = '2.976.96'
$ws.Range("E3").Value = '  +3.33%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '501.95'
$ws.Range("E5").Value = '  +5.05%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '135.02'
$ws.Range("E6").Value = '  +6.78%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.429'
$ws.Range("E8").Value = '  +6.20%  '
$ws.Range("E9").Value = '  +11.50%  '
$ws.Range("E10").Value = '  +10.60%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.349'
$ws.Range("E11").Value = '  +5.40%  '
$ws.Range("E12").Value = '  +3.65%  '
$ws.Range("D13").Value = '3.484.40'
$ws.Range("E13").Value = '  +3.36%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.40'
$ws.Range("E14").Value = '  +11.18%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000152'
$ws.Range("E15").Value = '  +12.50%  '
$ws.Range("D16").Value = '56.385.71'
$ws.Range("E16").Value = '  +3.48%  '
$ws.Range("D17").Value = '2.975.42'
$ws.Range("E17").Value = '  +3.20%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.74'
$ws.Range("E18").Value = '  +9.39%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.33'
$ws.Range("E19").Value = '  +6.96%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.74'
$ws.Range("E20").Value = '  +9.29%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '323.09'
$ws.Range("E21").Value = '  +5.58%  '
$ws.Range("E22").Value = '  +0.17%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.470'
$ws.Range("E23").Value = '  +5.48%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '61.72'
$ws.Range("E24").Value = '  +4.12%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.998'
$ws.Range("E25").Value = '  +0.02%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.163'
$ws.Range("E26").Value = '  +6.24%  '
$ws.Range("D27").Value = '0.0₃0891'
$ws.Range("E27").Value = '  +8.95%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.39'
$ws.Range("E28").Value = '  +1.45%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.77'
$ws.Range("E29").Value = '  +9.37%  '
$ws.Range("E30").Value = '  +4.27%  '
$ws.Range("E31").Value = '  +7.64%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.38'
$ws.Range("E32").Value = '  +6.89%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '157.87'
$ws.Range("E33").Value = '  +15.84%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.43'
$ws.Range("E34").Value = '  +4.98%  '
$ws.Range("E35").Value = '  +3.67%  '
$ws.Range("E36").Value = '  +1.54%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0669'
$ws.Range("E37").Value = '  +8.71%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '23.01'
$ws.Range("E38").Value = '  +0.94%  '
$ws.Range("D39").Value = '3.009.60'
$ws.Range("E39").Value = '  +3.61%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("E40").Value = '  -0.04%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.26'
$ws.Range("E41").Value = '  +1.82%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.641'
$ws.Range("E42").Value = '  +6.66%  '
$ws.Range("D43").Value = '2.243.05'
$ws.Range("E43").Value = '  +9.58%  '
$ws.Range("E44").Value = '  +5.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.978'
$ws.Range("E45").Value = '  +1.14%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.56'
$ws.Range("E46").Value = '  +4.15%  '
$ws.Range("E47").Value = '  +18.42%  '
$ws.Range("E48").Value = '  +9.77%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.73'
$ws.Range("E49").Value = '  +7.91%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.87'
$ws.Range("E50").Value = '  +5.57%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0865'
$ws.Range("E51").Value = '  +9.82%  '
